$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (2026-02-08 -> 2026-02-09) for every data row (rows 2 through 509).
$ws.Range("C2:C509").Value = 46062
